$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.059.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.808.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.10%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.19'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.47%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0683'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.070.03'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.849.24'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '11.04'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.65'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.011.67'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.68'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0788'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.43'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.72'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.46%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.24'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.85'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.88'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.48'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.56'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +17.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.11'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0556'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.03'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.76'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.701'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.50%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '92.04'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0193'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.312.40'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.51%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.43'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.24'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -7.04%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.17%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.96%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.987.08'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0666'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '99.52'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.57%  '
